# edit.ps1
# Applies the diff: rewrites several sentences in the Testing section,
# relocates the "Delivery" section to appear before "Design", updates
# "git" -> "GitHub" in the Design section, fixes "is was" -> "was" in the
# relocated Delivery text, and relocates the _GoBack bookmark to the end
# of the (now relocated) Delivery section.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) "One of our goals..." paragraph: trim the parenthetical aside and
#    reword "This is a new feature" -> "This introduced new functionality"
# ---------------------------------------------------------------
$find1 = "the beginning of the game (hereby referred to as “the new feature”). This is a new feature, so there we"
$repl1 = "the beginning of the game. This introduced new functionality, so there we"
$d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $repl1, 2) | Out-Null

# "A full set of unit tests were be created" -> "A full set of unit tests were created"
$find2 = "A full set of unit tests were be created"
$repl2 = "A full set of unit tests were created"
$d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $repl2, 2) | Out-Null

# " It has also" -> " This feature has also"
$find3 = " It has also"
$repl3 = " This feature has also"
$d.Content.Find.Execute($find3, $true, $false, $false, $false, $false, $true, 1, $false, $repl3, 2) | Out-Null

# ---------------------------------------------------------------
# 2) "For unit testing..." paragraph: drop the comparison-to-last-iteration
#    sentence.
# ---------------------------------------------------------------
$find4 = "% statement coverage for our non-graphics modules. This is approximately the same as the last iteration, however we plan to bring this number up some. "
$repl4 = "% statement coverage for our non-graphics modules."
$d.Content.Find.Execute($find4, $true, $false, $false, $false, $false, $true, 1, $false, $repl4, 2) | Out-Null

# ---------------------------------------------------------------
# 3) "The new feature is was developed" -> "The new feature was developed"
#    (fix the grammar slip left over from editing), done before the move
#    so the Find target text still matches verbatim.
# ---------------------------------------------------------------
$find5 = "The new feature is was developed"
$repl5 = "The new feature was developed"
$d.Content.Find.Execute($find5, $true, $false, $false, $false, $false, $true, 1, $false, $repl5, 2) | Out-Null

# ---------------------------------------------------------------
# 4) "git" -> "GitHub" in the Design section (whole-word match keeps the
#    existing run/proofErr wrapper around the word intact).
# ---------------------------------------------------------------
$find6 = "git"
$repl6 = "GitHub"
$d.Content.Find.Execute($find6, $true, $true, $false, $false, $false, $true, 1, $false, $repl6, 2) | Out-Null

# ---------------------------------------------------------------
# 5) Move the "Delivery" section (heading + its two paragraphs) so it
#    appears right before the "Design" heading instead of right after
#    the Design section's UML image.
# ---------------------------------------------------------------
$deliveryIdx = 0
$designIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($ptext -eq "Delivery" -and $deliveryIdx -eq 0) { $deliveryIdx = $i }
    if ($ptext -eq "Design" -and $designIdx -eq 0) { $designIdx = $i }
}

$srcStart = $d.Paragraphs.Item($deliveryIdx).Range.Start
$srcEnd = $d.Paragraphs.Item($deliveryIdx + 2).Range.End
$src = $d.Range($srcStart, $srcEnd)
$src.Cut() | Out-Null

$destPara = $d.Paragraphs.Item($designIdx)
$dest = $d.Range($destPara.Range.Start, $destPara.Range.Start)
$dest.Paste() | Out-Null

# ---------------------------------------------------------------
# 6) Relocate the "_GoBack" bookmark from around the UML picture to the
#    end of the (now relocated) Delivery section's last paragraph.
# ---------------------------------------------------------------
try {
    $oldBk = $d.Bookmarks.Item("_GoBack")
    $oldBk.Delete()
} catch {
}

$deliveryIdx2 = 0
$designIdx2 = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($ptext -eq "Delivery" -and $deliveryIdx2 -eq 0) { $deliveryIdx2 = $i }
    if ($ptext -eq "Design" -and $designIdx2 -eq 0) { $designIdx2 = $i }
}
$lastDeliveryPara = $d.Paragraphs.Item($deliveryIdx2 + 2)
$endPos = $lastDeliveryPara.Range.End - 2
$gbRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $gbRange) | Out-Null
